$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: add the dose (mg) and dose time for the second medication slot
$ws.Range("F14").Value = 20
$ws.Range("G14").Value = 0.56597222222222221
$ws.Range("G14").NumberFormat = "h:mm"

# Row 15: new day's entry - dose, dose time, BP time, and BP readings
# (only upload entries that have relevant data: bloodPressureDiastole,
# bloodPressureSystole, heartRate)
$ws.Range("C15").Value = 20
$ws.Range("D15").Value = 0.38750000000000001
$ws.Range("D15").NumberFormat = "h:mm"

$ws.Range("I15").Value = 0.38472222222222224
$ws.Range("I15").NumberFormat = "h:mm"

$ws.Range("J15").Value = 107
$ws.Range("K15").Value = 71
$ws.Range("L15").Value = 60

# Update the active selection (show graph data only for the current user)
$ws.Range("E15").Select()
